$wb = $excel.ActiveWorkbook

# Reference the original sheets by their existing names
$ws1 = $wb.Worksheets.Item("Shopenzer Testcases")
$ws2 = $wb.Worksheets.Item("Testscearnios")

# Rename sheets
$ws1.Name = "Car Resale Predictor Testcases"
$ws2.Name = "Test Scenarios"

# Make the second sheet ("Test Scenarios") the active/selected tab,
# with B22 as the active cell/selection
$ws2.Activate() | Out-Null
$ws2.Range("B22").Select() | Out-Null
